$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 84; this shifts all the existing rows (84..156) down
# by one (to 85..157), preserving their values/styles exactly as they were.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record's data.
$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value = "La Araucanía"
$ws.Cells.Item(84, 4).Value = 45271
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100101
$ws.Cells.Item(84, 8).Value = "Berries"
$ws.Cells.Item(84, 9).Value = 100101001
$ws.Cells.Item(84, 10).Value = "Arándano (blue)"
$ws.Cells.Item(84, 11).Value = "Sin especificar"
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 105
$ws.Cells.Item(84, 14).Value = 3000
$ws.Cells.Item(84, 15).Value = 3300
$ws.Cells.Item(84, 16).Value = 3157
$ws.Cells.Item(84, 17).Value = "`$/kilo"
$ws.Cells.Item(84, 18).Value = "Región del Maule"
$ws.Cells.Item(84, 19).Value = 3157
$ws.Cells.Item(84, 20).Value = 1
